$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "25.032.32"
$ws.Range("E2").Value = "  -3.65%  "
$ws.Range("D3").Value = "1.649.98"
$ws.Range("E3").Value = "  -5.44%  "
Set-TextValue $ws "D4" "1.000"
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws "D5" "237.16"
$ws.Range("E5").Value = "  -5.79%  "
Set-TextValue $ws "D6" "1.001"
$ws.Range("E6").Value = "  +0.09%  "
Set-TextValue $ws "D7" "0.4819"
$ws.Range("E7").Value = "  -6.28%  "
Set-TextValue $ws "D8" "0.2624"
$ws.Range("E8").Value = "  -5.01%  "
Set-TextValue $ws "D9" "0.06011"
$ws.Range("E9").Value = "  -2.98%  "
Set-TextValue $ws "D10" "0.07210"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").Value = "1.648.71"
$ws.Range("E11").Value = "  -5.50%  "
Set-TextValue $ws "D12" "14.86"
$ws.Range("E12").Value = "  -2.15%  "
Set-TextValue $ws "D13" "0.6220"
$ws.Range("E13").Value = "  -4.35%  "
Set-TextValue $ws "D14" "4.609"
$ws.Range("E14").Value = "  -0.59%  "
Set-TextValue $ws "D15" "73.02"
$ws.Range("E15").Value = "  -6.20%  "
Set-TextValue $ws "D16" "1.001"
$ws.Range("E16").Value = "  +0.08%  "
Set-TextValue $ws "D17" "1.000"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "25.019.29"
$ws.Range("E18").Value = "  -3.80%  "
$ws.Range("E19").Value = "  -2.75%  "
Set-TextValue $ws "D20" "0.000006632"
$ws.Range("E20").Value = "  -2.59%  "
Set-TextValue $ws "D21" "4.588"
$ws.Range("E21").Value = "  +6.48%  "
$ws.Range("D22").Value = "1.856.56"
$ws.Range("E22").Value = "  -5.60%  "
Set-TextValue $ws "D23" "8.638"
$ws.Range("E23").Value = "  -0.63%  "
Set-TextValue $ws "D24" "5.307"
$ws.Range("E24").Value = "  -1.43%  "
Set-TextValue $ws "D25" "132.26"
$ws.Range("E25").Value = "  -2.68%  "
Set-TextValue $ws "D26" "14.97"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("E27").Value = "  -7.72%  "
Set-TextValue $ws "D28" "103.49"
$ws.Range("E28").Value = "  -2.49%  "
Set-TextValue $ws "D29" "1.671"
$ws.Range("E29").Value = "  -6.48%  "
Set-TextValue $ws "D30" "3.777"
$ws.Range("E30").Value = "  -4.49%  "
Set-TextValue $ws "D31" "0.07905"
$ws.Range("E31").Value = "  -4.22%  "
Set-TextValue $ws "D32" "3.607"
$ws.Range("E32").Value = "  -2.00%  "
Set-TextValue $ws "D33" "0.04625"
$ws.Range("E33").Value = "  -1.31%  "
Set-TextValue $ws "D34" "2.597"
$ws.Range("E34").Value = "  -2.16%  "
Set-TextValue $ws "D35" "0.9396"
$ws.Range("E35").Value = "  -6.07%  "
Set-TextValue $ws "D36" "0.5781"
$ws.Range("E36").Value = "  -7.68%  "
Set-TextValue $ws "D37" "2.604"
$ws.Range("E37").Value = "  -4.88%  "
Set-TextValue $ws "D38" "0.01565"
$ws.Range("E38").Value = "  -3.01%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D39" "0.8410"
$ws.Range("E39").Value = "  +10.00%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws "D40" "1.000"
$ws.Range("E40").Value = "  +0.11%  "
Set-TextValue $ws "D41" "1.831"
$ws.Range("E41").Value = "  -4.70%  "
Set-TextValue $ws "D42" "98.28"
$ws.Range("E42").Value = "  -2.36%  "
Set-TextValue $ws "D43" "0.3739"
$ws.Range("E43").Value = "  -3.85%  "
Set-TextValue $ws "D44" "4.792"
$ws.Range("E44").Value = "  -4.68%  "
Set-TextValue $ws "D45" "0.1143"
$ws.Range("E45").Value = "  +0.90%  "
Set-TextValue $ws "D46" "6.165"
$ws.Range("E46").Value = "  -3.00%  "
Set-TextValue $ws "D47" "0.05194"
$ws.Range("E47").Value = "  -0.71%  "
Set-TextValue $ws "D48" "29.88"
$ws.Range("E48").Value = "  -3.01%  "
Set-TextValue $ws "D49" "50.77"
$ws.Range("E49").Value = "  -8.66%  "
$ws.Range("E50").Value = "  -0.07%  "
Set-TextValue $ws "D51" "0.3344"
$ws.Range("E51").Value = "  -3.20%  "
